# This script reorders the "Criterion N, Air Speed 0.1" worksheets so that the
# data/labels rotate: sheet that used to be "Criterion 1" becomes "Criterion 2",
# the one that used to be "Criterion 3" becomes "Criterion 1", and the one that
# used to be "Criterion 2" becomes "Criterion 3". It also reorders the columns
# on the readme sheet/table (Author moves before JobNo/sheet_name).

$wb = $excel.ActiveWorkbook

# --- Grab worksheet references by their current (pre-edit) names -----------
$wsReadme = $wb.Worksheets.Item("readme")
$wsCrit1  = $wb.Worksheets.Item("Criterion 1, Air Speed 0.1")
$wsCrit3  = $wb.Worksheets.Item("Criterion 3, Air Speed 0.1")
$wsCrit2  = $wb.Worksheets.Item("Criterion 2, Air Speed 0.1")

# --- Read the full data grids (rows 2-32, columns C:F) before changing -----
# anything, so none of the source values gets clobbered mid-way.
$gridCrit1 = @{}
$gridCrit3 = @{}
$gridCrit2 = @{}

for ($r = 2; $r -le 32; $r++) {
    for ($c = 3; $c -le 6; $c++) {
        $key = "$r,$c"
        $gridCrit1[$key] = $wsCrit1.Cells.Item($r, $c).Value()
        $gridCrit3[$key] = $wsCrit3.Cells.Item($r, $c).Value()
        $gridCrit2[$key] = $wsCrit2.Cells.Item($r, $c).Value()
    }
}

# --- Read readme sheet values (rows 2-5) before changing columns -----------
# (column E "Date" is not touched at all, so we don't need to read/rewrite it)
$readmeJobNo     = @{}
$readmeSheetName = @{}
$readmeAuthor    = @{}
for ($r = 2; $r -le 5; $r++) {
    $readmeJobNo[$r]     = $wsReadme.Cells.Item($r, 2).Value()
    $readmeSheetName[$r] = $wsReadme.Cells.Item($r, 3).Value()
    $readmeAuthor[$r]    = $wsReadme.Cells.Item($r, 4).Value()
}

function Write-CritCell($ws, $r, $c, $val) {
    if ($val -eq $null) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $val
    }
}

# --- Write rotated data into the destination sheets -------------------------
# wsCrit1 (tab "Criterion 1, Air Speed 0.1") becomes "Criterion 2" -> gets old Criterion 2 data
# wsCrit3 (tab "Criterion 3, Air Speed 0.1") becomes "Criterion 1" -> gets old Criterion 1 data
# wsCrit2 (tab "Criterion 2, Air Speed 0.1") becomes "Criterion 3" -> gets old Criterion 3 data
for ($r = 2; $r -le 32; $r++) {
    for ($c = 3; $c -le 6; $c++) {
        $key = "$r,$c"
        Write-CritCell $wsCrit1 $r $c $gridCrit2[$key]
        Write-CritCell $wsCrit3 $r $c $gridCrit1[$key]
        Write-CritCell $wsCrit2 $r $c $gridCrit3[$key]
    }
}

# --- Update header labels (row 1, columns E/F) for each sheet's new identity
$wsCrit1.Cells.Item(1, 5).Value = "Criterion 2 Absolute Change"
$wsCrit1.Cells.Item(1, 6).Value = "Criterion 2 Relative Change (%)"

$wsCrit3.Cells.Item(1, 5).Value = "Criterion 1 Absolute Change"
$wsCrit3.Cells.Item(1, 6).Value = "Criterion 1 Relative Change (%)"

$wsCrit2.Cells.Item(1, 5).Value = "Criterion 3 Absolute Change"
$wsCrit2.Cells.Item(1, 6).Value = "Criterion 3 Relative Change (%)"

# --- Rename the worksheet tabs to match their new identity ------------------
# (use temporary names first since this is a cyclic rename and the target
# names are currently in use by other sheets)
$wsCrit1.Name = "__tmp_Crit1__"
$wsCrit3.Name = "__tmp_Crit3__"
$wsCrit2.Name = "__tmp_Crit2__"

$wsCrit1.Name = "Criterion 2, Air Speed 0.1"
$wsCrit3.Name = "Criterion 1, Air Speed 0.1"
$wsCrit2.Name = "Criterion 3, Air Speed 0.1"

# --- Update the readme sheet: reorder columns to index,Author,JobNo,sheet_name,Date
$wsReadme.Cells.Item(1, 2).Value = "Author"
$wsReadme.Cells.Item(1, 3).Value = "JobNo"
$wsReadme.Cells.Item(1, 4).Value = "sheet_name"
$wsReadme.Cells.Item(1, 5).Value = "Date"

# sheet_name values rotate the same way the tabs were renamed above
$sheetNameRotation = @{
    "Criteria Failing, Air Speed 0.1" = "Criteria Failing, Air Speed 0.1"
    "Criterion 1, Air Speed 0.1"      = "Criterion 2, Air Speed 0.1"
    "Criterion 3, Air Speed 0.1"      = "Criterion 1, Air Speed 0.1"
    "Criterion 2, Air Speed 0.1"      = "Criterion 3, Air Speed 0.1"
}

for ($r = 2; $r -le 5; $r++) {
    $oldSheetName = $readmeSheetName[$r]
    $newSheetName = $sheetNameRotation[$oldSheetName]
    $wsReadme.Cells.Item($r, 2).Value = $readmeAuthor[$r]
    $wsReadme.Cells.Item($r, 3).Value = $readmeJobNo[$r]
    $wsReadme.Cells.Item($r, 4).Value = $newSheetName
}
